$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-15
# from 2023-10-09 (45208) to 2023-10-13 (45212)
$ws.Range("C2:C15").Value = 45212
